$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("work")

# --- Row 24 (job on 43912) was left open; fill in its FINISHED time and HOURS ---
# Copy the number formatting from the row above (which already has all six columns
# filled in) onto the two still-blank cells before writing their values.
$ws.Range("E23:F23").Copy()
$ws.Range("E24:F24").PasteSpecial(-4122)  # xlPasteFormats

$ws.Cells.Item(24, 5).Value = 0.14583333333333334   # FINISHED 3:30 AM
$ws.Cells.Item(24, 6).Value = 2                      # HOURS

# --- New row 25: a second job entry ("Merge files") for the same day ---
$ws.Range("B24:D24").Copy()
$ws.Range("B25:D25").PasteSpecial(-4122)  # xlPasteFormats

$ws.Cells.Item(25, 2).Value = "2"
$ws.Cells.Item(25, 3).Value = "Merge files"
$ws.Cells.Item(25, 4).Value = 0.14583333333333334   # BEGIN 3:30 AM

$ws.Application.CutCopyMode = $false

# --- View state: selection moves to C26 (first empty row under the new data) ---
$ws.Range("C26").Select()
